$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# The Price column stores plain-text values (e.g. "1.008", "20.91") in the
# original workbook. Assigning such numeric-looking strings directly to
# .Value would make Excel auto-convert them into real numbers, so each
# cell is temporarily switched to a Text number format before the value is
# written, then restored to the default (Normal) style afterwards so the
# cell keeps its original (unstyled) appearance.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.310.88'
$ws.Range("D3").Value = '1.666.86'
$ws.Range("D4").Value = '1.008'
$ws.Range("D6").Value = '0.5307'
$ws.Range("D7").Value = '1.009'
$ws.Range("D9").Value = '0.06358'
$ws.Range("D10").Value = '20.91'
$ws.Range("D11").Value = '0.07836'
$ws.Range("D12").Value = '4.529'
$ws.Range("D13").Value = '1.669.12'
$ws.Range("D14").Value = '1.895.50'
$ws.Range("D15").Value = '0.5605'
$ws.Range("D16").Value = '0.0₅8129'
$ws.Range("D17").Value = '65.75'
$ws.Range("D18").Value = '26.318.57'
$ws.Range("D21").Value = '199.04'
$ws.Range("D22").Value = '10.27'
$ws.Range("D23").Value = '6.051'
$ws.Range("D25").Value = '146.62'
$ws.Range("D26").Value = '0.1213'
$ws.Range("D28").Value = '16.14'
$ws.Range("D29").Value = '1.512'
$ws.Range("D30").Value = '0.05887'
$ws.Range("D32").Value = '3.536'
$ws.Range("D34").Value = '1.601'
$ws.Range("D35").Value = '2.830'
$ws.Range("D36").Value = '0.9600'
$ws.Range("D38").Value = '0.5801'
$ws.Range("D39").Value = '0.01618'
$ws.Range("D41").Value = '1.073.97'
$ws.Range("D42").Value = '0.8568'
$ws.Range("D44").Value = '102.76'
$ws.Range("D45").Value = '1.806.27'
$ws.Range("D46").Value = '58.40'
$ws.Range("D48").Value = '0.4414'
$ws.Range("D50").Value = '8.038'
$ws.Range("D51").Value = '0.05147'

$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E) updates ---
# These values include a leading "  " and trailing "  " plus a "%" sign, so
# Excel keeps them as text automatically - no special handling required.
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +1.67%  '
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("E21").Value = '  +3.68%  '
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("E30").Value = '  +1.73%  '
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  -0.65%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +2.08%  '
$ws.Range("E41").Value = '  +2.79%  '
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("E51").Value = '  -0.25%  '
